# estadisticas_2026.xlsx — "Add files via upload" update
#
# Eventos sheet:
#  - row 36 (player 21): a yellow-card/red-card event is now recorded
#    (gol_recibido -> 1, rojas -> 1)
#  - a brand new event row is inserted right before the last row (old row 46,
#    player 13) for player 53 on team "azul" with gol_recibido = 7; the
#    previously-last row (player 13) shifts down to row 47 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

# --- row 36: mark gol_recibido / rojas -------------------------------------
$ws.Range("D36").Value = 1
$ws.Range("L36").Value = 1

# --- insert a fresh row before the current last row (row 46) ---------------
# This pushes the existing row 46 (and its formatting / shared H formula)
# down to row 47, and leaves a blank row 46 behind for the new record.
$ws.Rows.Item(46).Insert()

# --- populate the newly inserted row 46 with the new event record ----------
$ws.Range("A46").Value = 2
$ws.Range("B46").Value = 53
$ws.Range("C46").Value = "azul"
$ws.Range("D46").Value = 7
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0

# --- move the on-screen selection to where the editor left it --------------
$ws.Range("N46").Select()

# --- best-effort: scroll the Jugadores sheet so row 35 is near the top -----
# (kept last / wrapped in try so it can never abort the rest of the script)
try {
    $wsJugadores = $wb.Worksheets.Item("Jugadores")
    $wsJugadores.Activate()
    $excel.ActiveWindow.ScrollRow = 35
} catch {
} finally {
    $ws.Activate()
}
